$d = $word.ActiveDocument

# --- Fill in the three "0€" placeholder cells on the first data row of the
#     expense table (row 2 of Table 1: Repas / Hebergement / Total columns) ---
$table = $d.Tables.Item(1)
$table.Cell(2, 7).Range.Text = "999€"
$table.Cell(2, 8).Range.Text = "222€"
$table.Cell(2, 9).Range.Text = "1€"

# --- Update the total amount (7440,00 -> 7441,00) everywhere it appears ---
$d.Content.Find.Execute("7440,00", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "7441,00", 2)
